$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.299009561538696
$ws.Range("B1").Value = 1.80309784412384
$ws.Range("C1").Value = 4.235278606414795
$ws.Range("D1").Value = 2.977567911148071
$ws.Range("E1").Value = 1.134094715118408
